$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 697
$ws.Range("I2").Value  = 1859
$ws.Range("J2").Value  = 7458
$ws.Range("K2").Value  = 35
$ws.Range("L2").Value  = 2153
$ws.Range("M2").Value  = 123
$ws.Range("N2").Value  = 1316
$ws.Range("O2").Value  = 7
$ws.Range("P2").Value  = 33
$ws.Range("Q2").Value  = 17
$ws.Range("R2").Value  = 123
$ws.Range("S2").Value  = 799
$ws.Range("T2").Value  = 1281
$ws.Range("U2").Value  = 100
$ws.Range("V2").Value  = 11632
$ws.Range("W2").Value  = 2
$ws.Range("X2").Value  = 11570
$ws.Range("Y2").Value  = 13
$ws.Range("Z2").Value  = 192
$ws.Range("AA2").Value = 69
